$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'TOLENTINO VASQUEZ DIANA KATHERYN'
$ws.Cells.Item(2, 2).Value = 70

$ws.Cells.Item(3, 1).Value = 'JULCA VALENZUELA CINTIA KARYN'
$ws.Cells.Item(3, 2).Value = 68

$ws.Cells.Item(4, 1).Value = 'SANCHEZ CORTEZ LEYLA DIANA'
$ws.Cells.Item(4, 2).Value = 64

$ws.Cells.Item(5, 1).Value = 'VALER VEGA PATRICIA GERALDINE'
$ws.Cells.Item(5, 2).Value = 61

$ws.Cells.Item(6, 1).Value = 'YZQUIERDO CARHUATANTA LEYDY YANELA'
$ws.Cells.Item(6, 2).Value = 58

$ws.Cells.Item(7, 1).Value = 'RODRIGUEZ RUBIO SANDRA MABEL'
$ws.Cells.Item(7, 2).Value = 58

$ws.Cells.Item(8, 1).Value = 'DE LA CRUZ BENITES RICHARD ALEXANDER'
$ws.Cells.Item(8, 2).Value = 58

$ws.Cells.Item(9, 1).Value = 'PONCE VILLANUEVA CARMEN ISABEL'
$ws.Cells.Item(9, 2).Value = 56

$ws.Cells.Item(10, 1).Value = 'ARENAS ZAVALA ANDYELA PATRICIA ISIDORA'
$ws.Cells.Item(10, 2).Value = 55

$ws.Cells.Item(11, 1).Value = 'GASLAC GUTIERREZ FRANK JHORDY'
$ws.Cells.Item(11, 2).Value = 53

$ws.Cells.Item(12, 1).Value = 'RUBIO MARIÑOS GISELA JUDITH'
$ws.Cells.Item(12, 2).Value = 52

$ws.Cells.Item(13, 1).Value = 'REYES RODRIGUEZ JEISSON STEVEN'
$ws.Cells.Item(13, 2).Value = 52

$ws.Cells.Item(14, 1).Value = 'CARRILLO MARTÍNEZ HEIDY NAYELI'
$ws.Cells.Item(14, 2).Value = 51

$ws.Cells.Item(15, 1).Value = 'CYNTHIA RODRIGUEZ LECCA'
$ws.Cells.Item(15, 2).Value = 50

$ws.Cells.Item(16, 1).Value = 'PIERINA NAGIELLY SANDOVAL CONTRERAS'
$ws.Cells.Item(16, 2).Value = 49

$ws.Cells.Item(17, 1).Value = 'SEGURA ASTO YAMILET ANTONELA'
$ws.Cells.Item(17, 2).Value = 49

$ws.Cells.Item(18, 1).Value = 'GUZMAN ZAVALETA CECILIA MARISOL'
$ws.Cells.Item(18, 2).Value = 43

$ws.Cells.Item(19, 1).Value = 'LEON VERA MELISSA FIORELLA'
$ws.Cells.Item(19, 2).Value = 30

$ws.Cells.Item(20, 1).Value = 'RODRIGUEZ VASQUEZ WALTER'
$ws.Cells.Item(20, 2).Value = 30

$ws.Cells.Item(21, 1).Value = 'RUTH MELISSA RAMIREZ VELEZMORO'
$ws.Cells.Item(21, 2).Value = 16
